$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: update the payment-approval record for 0420172008483 ---
# Set E5 (NroSiniestro) before C5 (Usuario) so the shared-string table is
# built in the same order the original authoring session produced.
$ws.Range("E5").Value = "'0420172008483   "
$ws.Range("C5").Value = "nbedoya"

# --- Row 6: new record, cloned from row 5's layout/format ---
$ws.Range("A5:E5").Copy($ws.Range("A6:E6"))
$ws.Range("C6").Value = "rsuarez"
$ws.Range("E6").Value = "'1220170301429   "
$ws.Hyperlinks.Add($ws.Range("B6"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B6").Style = "Hipervínculo"

# --- Row 7: another new record, cloned from row 5's layout/format ---
$ws.Range("A5:E5").Copy($ws.Range("A7:E7"))
$ws.Range("C7").Value = "dgariffo"
$ws.Range("E7").Value = "'1120170200936   "
$ws.Hyperlinks.Add($ws.Range("B7"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B7").Style = "Hipervínculo"

# Leave the selection where the author left it after entering the new rows.
$ws.Range("H8").Select()
